# Add two new task rows to the "윤다은" sheet's Table1 (A1:F17 -> A1:F19).
#
# Edit history being replicated (inferred from the shared-string order in the
# target file): on 2019-06-18 a "testing" task was logged as the last row of
# the table; the next day (2019-06-19) a new "예약 업로드" task was inserted
# above it, pushing the "testing" row down by one. We reproduce that by
# writing the "testing" row's cells (row 19) before the "예약 업로드" row's
# cells (row 18), which yields the same shared-string allocation order as
# the target workbook, even though row 18 ends up above row 19 on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("윤다은")
$lo = $ws.ListObjects.Item("Table1")

# Bring formatting (wrap/vertical-top text style for A,B,E,F and the custom
# date format for C,D) down into the two new rows by copying the last
# existing data row's formatting.
$ws.Range("A17:F17").Copy() | Out-Null
$ws.Range("A18:F19").PasteSpecial(-4122) | Out-Null

# Row 19 ("testing" / 2019-06-18) — written first.
$ws.Cells.Item(19, 1).Value = "testing"
$ws.Cells.Item(19, 2).Value = "에러탐지 및 줄이기"
$ws.Cells.Item(19, 3).Value = 43634
$ws.Cells.Item(19, 5).Value = "new posting 할 때, 침대수와 인원수, 가격은 숫자만 들어올 수 있도록 처리//"
$ws.Rows(19).RowHeight = 51.75

# Row 18 ("예약 업로드" / 2019-06-19) — written second.
$ws.Cells.Item(18, 1).Value = "예약 업로드"
$ws.Cells.Item(18, 2).Value = "날짜 분석해서 겹치지 않게 디비에 업로드"
$ws.Cells.Item(18, 3).Value = 43635

# Grow the table (and its AutoFilter range) to cover the two new rows.
$lo.Resize($ws.Range("A1:F19"))

# Match the saved selection/scroll position recorded in the target file.
$ws.Activate()
$ws.Range("D18").Select() | Out-Null
